# Applies the Alvearie FHIR IG metadata refresh to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (first sheet) ---
$ws1 = $wb.Worksheets.Item(1)

# Version bump 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date refresh
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value
$ws1.Range("B9").Value = "Alvearie Team"

# Replace the duplicated "Contact" row (row 10) with a new "Jurisdiction" row
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Remove the now-redundant duplicate "Contact" row (old row 11), shifting rows 12-21 up
$ws1.Rows.Item(11).Delete()

# --- Sheet "Elements" (second sheet) ---
$ws2 = $wb.Worksheets.Item(2)

# Update the Short/Definition text for the root Extension slice
$ws2.Range("K2").Value = "Immigration Status"
$ws2.Range("L2").Value = "Customer-specific code for the immigration status of the person"
